# Update scripts with new TPM values (NATMI LR-pairs output: Col2a1-Mag)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.06437833333333333
$ws.Range("H2").Value = 0.193135
$ws.Range("I2").Value = 0.109187438766332
$ws.Range("J2").Value = 0.109187438766332
$ws.Range("M2").Value = 0.232947
$ws.Range("N2").Value = 0.698841
$ws.Range("O2").Value = 0.2572219815457369
$ws.Range("P2").Value = 0.2572219815457369
$ws.Range("Q2").Value = 0.014996739615
$ws.Range("R2").Value = 0.134970656535
$ws.Range("S2").Value = 0.02808540935937974
$ws.Range("T2").Value = 0.02808540935937974

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.06437833333333333
$ws.Range("H3").Value = 0.193135
$ws.Range("I3").Value = 0.109187438766332
$ws.Range("J3").Value = 0.109187438766332
$ws.Range("M3").Value = 0.6726793333333333
$ws.Range("N3").Value = 2.018038
$ws.Range("O3").Value = 0.7427780184542632
$ws.Range("P3").Value = 0.7427780184542632
$ws.Range("Q3").Value = 0.04330597434777778
$ws.Range("R3").Value = 0.38975376913
$ws.Range("S3").Value = 0.08110202940695231
$ws.Range("T3").Value = 0.08110202940695231

# Row 4
$ws.Range("G4").Value = 0.4788196666666666
$ws.Range("I4").Value = 0.8120914339857952
$ws.Range("J4").Value = 0.8120914339857951
$ws.Range("M4").Value = 0.232947
$ws.Range("N4").Value = 0.698841
$ws.Range("O4").Value = 0.2572219815457369
$ws.Range("P4").Value = 0.2572219815457369
$ws.Range("Q4").Value = 0.111539604891
$ws.Range("R4").Value = 1.003856444019
$ws.Range("S4").Value = 0.2088877678461452
$ws.Range("T4").Value = 0.2088877678461452

# Row 5
$ws.Range("G5").Value = 0.4788196666666666
$ws.Range("I5").Value = 0.8120914339857952
$ws.Range("J5").Value = 0.8120914339857951
$ws.Range("M5").Value = 0.6726793333333333
$ws.Range("N5").Value = 2.018038
$ws.Range("O5").Value = 0.7427780184542632
$ws.Range("P5").Value = 0.7427780184542632
$ws.Range("Q5").Value = 0.3220920941602222
$ws.Range("R5").Value = 2.898828847442
$ws.Range("S5").Value = 0.60320366613965
$ws.Range("T5").Value = 0.6032036661396499

# Row 6
$ws.Range("G6").Value = 0.042481
$ws.Range("H6").Value = 0.127443
$ws.Range("I6").Value = 0.07204895414449818
$ws.Range("J6").Value = 0.07204895414449818
$ws.Range("M6").Value = 0.232947
$ws.Range("N6").Value = 0.698841
$ws.Range("O6").Value = 0.2572219815457369
$ws.Range("P6").Value = 0.2572219815457369
$ws.Range("Q6").Value = 0.009895821507
$ws.Range("R6").Value = 0.08906239356300001
$ws.Range("S6").Value = 0.01853257475334575
$ws.Range("T6").Value = 0.01853257475334575

# Row 7
$ws.Range("G7").Value = 0.042481
$ws.Range("H7").Value = 0.127443
$ws.Range("I7").Value = 0.07204895414449818
$ws.Range("J7").Value = 0.07204895414449818
$ws.Range("M7").Value = 0.6726793333333333
$ws.Range("N7").Value = 2.018038
$ws.Range("O7").Value = 0.7427780184542632
$ws.Range("P7").Value = 0.7427780184542632
$ws.Range("Q7").Value = 0.02857609075933333
$ws.Range("R7").Value = 0.257184816834
$ws.Range("S7").Value = 0.05351637939115243
$ws.Range("T7").Value = 0.05351637939115243

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.003934
$ws.Range("H8").Value = 0.011802
$ws.Range("I8").Value = 0.006672173103374587
$ws.Range("J8").Value = 0.006672173103374586
$ws.Range("M8").Value = 0.232947
$ws.Range("N8").Value = 0.698841
$ws.Range("O8").Value = 0.2572219815457369
$ws.Range("P8").Value = 0.2572219815457369
$ws.Range("Q8").Value = 0.0009164134980000001
$ws.Range("R8").Value = 0.008247721482
$ws.Range("S8").Value = 0.00171622958686618
$ws.Range("T8").Value = 0.00171622958686618

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.003934
$ws.Range("H9").Value = 0.011802
$ws.Range("I9").Value = 0.006672173103374587
$ws.Range("J9").Value = 0.006672173103374586
$ws.Range("M9").Value = 0.6726793333333333
$ws.Range("N9").Value = 2.018038
$ws.Range("O9").Value = 0.7427780184542632
$ws.Range("P9").Value = 0.7427780184542632
$ws.Range("Q9").Value = 0.002646320497333333
$ws.Range("R9").Value = 0.023816884476
$ws.Range("S9").Value = 0.004955943516508407
$ws.Range("T9").Value = 0.004955943516508406
